# "Generate Report for Handback"
#
# The CI localization-status report is refreshed once handback data is
# available: the Overview sheet's status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", and each language sheet (zh-cn, de-de)
# gets its "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns filled in (plus a hyperlink on the newly-populated
# "Latest Target File" cell, matching the existing hyperlink style already
# used elsewhere in the workbook).

$wb = $excel.ActiveWorkbook

$mdFileName = "1cb91d82-c797-4adc-bdaf-aeab64b649ef.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/def91282140063ca3cb0f64079699953de104eec/e2e/1cb91d82-c797-4adc-bdaf-aeab64b649ef.md"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status column for both languages
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Widen the status columns to fit the longer text (matches the other
# "wide" columns already present in the workbook).
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus

$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Range("J2").Value = "1cb91d82-c797-4adc-bdaf-aeab64b649ef.f90b542ae74e50539d3d5afef1203f0252100bff.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-22 15:08:44"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null

$wsZhCn.Columns.Item(3).ColumnWidth = 29.15
$wsZhCn.Columns.Item(9).ColumnWidth = 39.15
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus

$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Range("J2").Value = "1cb91d82-c797-4adc-bdaf-aeab64b649ef.f90b542ae74e50539d3d5afef1203f0252100bff.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-22 15:08:51"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null

$wsDeDe.Columns.Item(3).ColumnWidth = 29.15
$wsDeDe.Columns.Item(9).ColumnWidth = 39.15
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15
